# movimientos.xlsx - update the single data row (row 2) with the latest
# balances/totals. Values are stored as literal text (matching the sheet's
# existing convention of thousands-separated strings rather than numeric
# cells), so a leading apostrophe forces Excel to keep them as text instead
# of auto-converting the comma-formatted numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'.00"
$ws.Range("B2").Value = "'545,474,228.14"
$ws.Range("C2").Value = "'1,570.00"
$ws.Range("D2").Value = "'545,472,658.14"
$ws.Range("E2").Value = "'22,416.14"
